$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value (45204 -> 2023-10-05).
# Update every data row (2 through 146) to the new date serial 45207 (2023-10-08).
$ws.Range("C2:C146").Value = 45207
